$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 219-220; this shifts the existing rows 219..264
# down to 221..266 (and auto-extends the sheet dimension to A1:R266).
$ws.Rows("219:220").Insert()

# Row 219 (new): same Mercado/Region/Categoria/Variedad/Calidad as the
# original row 219, but with an updated date, volume, price range/avg,
# unit size, origin and derived per-kg price.
$ws.Range("A219").Value2 = 9
$ws.Range("B219").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C219").Value2 = "Metropolitana"
$ws.Range("D219").Value2 = 44637
$ws.Range("E219").Value2 = 13
$ws.Range("F219").Value2 = 100112021
$ws.Range("G219").Value2 = "Ají"
$ws.Range("H219").Value2 = "Inferno"
$ws.Range("I219").Value2 = "Primera"
$ws.Range("J219").Value2 = 65
$ws.Range("K219").Value2 = 20000
$ws.Range("L219").Value2 = 20000
$ws.Range("M219").Value2 = 20000
$ws.Range("N219").Value2 = "$/caja 15 kilos"
$ws.Range("O219").Value2 = "Provincia de Quillota"
$ws.Range("P219").Value2 = 1333
$ws.Range("Q219").Value2 = 15
$ws.Range("R219").Value2 = "Hortaliza"

# Row 220 (new): a second new observation for the same date.
$ws.Range("A220").Value2 = 9
$ws.Range("B220").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C220").Value2 = "Metropolitana"
$ws.Range("D220").Value2 = 44637
$ws.Range("E220").Value2 = 13
$ws.Range("F220").Value2 = 100112021
$ws.Range("G220").Value2 = "Ají"
$ws.Range("H220").Value2 = "Inferno"
$ws.Range("I220").Value2 = "Segunda"
$ws.Range("J220").Value2 = 42
$ws.Range("K220").Value2 = 16000
$ws.Range("L220").Value2 = 16000
$ws.Range("M220").Value2 = 16000
$ws.Range("N220").Value2 = "$/caja 15 kilos"
$ws.Range("O220").Value2 = "Provincia de Quillota"
$ws.Range("P220").Value2 = 1067
$ws.Range("Q220").Value2 = 15
$ws.Range("R220").Value2 = "Hortaliza"
